# #5: property boat&car done
# Sheet "汽車" (car/boat property) row 1 was a stray duplicate of the data
# row instead of real column headers; fix it to proper headers and extend
# the table with the metadata columns (property_category .. index) plus a
# "capacity" column, matching the other property sheets in this workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# ---- Row 1: turn the accidental duplicate data row into real headers ----
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"

$newHeaderCols = @("H", "I", "J", "K", "L", "M", "N")
$newHeaderVals = @(
    "property_category",
    "category",
    "date",
    "legislator_name",
    "legislator_id",
    "source_file",
    "index"
)

for ($i = 0; $i -lt $newHeaderCols.Length; $i++) {
    $cell = $ws.Range($newHeaderCols[$i] + "1")
    $cell.Value = $newHeaderVals[$i]
}

# B1:G1 already carry the bold/centered/bordered header look; only the
# newly-added H1:N1 cells need it applied so they match.
$newHeaderRange = $ws.Range("H1:N1")
$newHeaderRange.Borders.LineStyle = 1
$newHeaderRange.HorizontalAlignment = -4108
$newHeaderRange.VerticalAlignment = -4160
$newHeaderRange.Font.Bold = $true

# ---- Row 2: keep existing data (A2:G2), append the new tracked columns ----
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"

# "2012-04-30" looks like a date, so COM auto-parses it into a date serial
# unless the cell is forced to text first; clear the forced format again
# afterwards so it doesn't leave a date-format override behind.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2012-04-30"
$ws.Range("J2").ClearFormats()

$ws.Range("K2").Value = "姚文智"
$ws.Range("L2").Value = 1745
$ws.Range("M2").Value = "tmp60da1"
$ws.Range("N2").Value = 28
